$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "pp00065"
$ws.Range("B4").Value = "3lü Priz"
$ws.Range("B5").Value = "Kum"
$ws.Range("B6").Value = "Dökmelik"
$ws.Range("B7").Value = ""

$ws.Range("B9").Value = 43
$ws.Range("B10").Value = 6
$ws.Range("B11").Value = 3
$ws.Range("B13").Value = 65
$ws.Range("B14").Value = 47
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "6"
$ws.Range("B16").Value = 89
$ws.Range("B17").Value = 68
$ws.Range("B18").Value = 6
$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 9
$ws.Range("B21").Value = 7
$ws.Range("B22").Value = 0
$ws.Range("B23").Value = 43
$ws.Range("B24").Value = 0
$ws.Range("B25").Value = 0
$ws.Range("B26").Value = 0
$ws.Range("B27").Value = 0
$ws.Range("B28").Value = 0
$ws.Range("B29").Value = 0
$ws.Range("B30").Value = 0
$ws.Range("B31").Value = 0
$ws.Range("B32").Value = 0
$ws.Range("B33").Value = 0
$ws.Range("B34").Value = 0
$ws.Range("B35").Value = 0
$ws.Range("B36").Value = 0
$ws.Range("B37").Value = 0
$ws.Range("B38").Value = 0
$ws.Range("B39").Value = 0
$ws.Range("B40").Value = 0
$ws.Range("B41").Value = 0
$ws.Range("B42").Value = 0
$ws.Range("B43").Value = 0
$ws.Range("B44").Value = 0
$ws.Range("B45").Value = 0
$ws.Range("B46").Value = 0
$ws.Range("B47").Value = 0
